$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1093.6666
$ws.Range("I43").Value = 870
$ws.Range("J43").Value = 1138.4
$ws.Range("K43").Value = 870
$ws.Range("L43").Value = 1138.4
$ws.Range("M43").Value = -801
$ws.Range("N43").Value = -1276.4
$ws.Range("H69").Value = 20963928
$ws.Range("I69").Value = 2802.1667
$ws.Range("J69").Value = 23697986
$ws.Range("K69").Value = 8406.500100000001
$ws.Range("L69").Value = 71093958
$ws.Range("M69").Value = -7532.500100000001
$ws.Range("N69").Value = -71095706
$ws.Range("H72").Value = 20963928
$ws.Range("I72").Value = 2802.1667
$ws.Range("J72").Value = 23697986
$ws.Range("K72").Value = 25219.5003
$ws.Range("L72").Value = 213281874
$ws.Range("M72").Value = -20851.5003
$ws.Range("N72").Value = -213290610
$ws.Range("H116").Value = 3347011.2
$ws.Range("I116").Value = 12822925
$ws.Range("J116").Value = 2570.9412
$ws.Range("K116").Value = 12822925
$ws.Range("L116").Value = 2570.9412
$ws.Range("M116").Value = -12819483
$ws.Range("N116").Value = -9454.941200000001
$ws.Range("H134").Value = 37946.316
$ws.Range("J134").Value = 37946.316
$ws.Range("L134").Value = 37946.316
$ws.Range("N134").Value = -48086.316
$ws.Range("H138").Value = 3601.1094
$ws.Range("I138").Value = 1341.5714
$ws.Range("J138").Value = 4233.78
$ws.Range("K138").Value = 4024.7142
$ws.Range("L138").Value = 12701.34
$ws.Range("M138").Value = 1115.2858
$ws.Range("N138").Value = -22981.34

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 32980
$ws.Range("J138").Value = 32980
$ws.Range("L138").Value = 32980
$ws.Range("N138").Value = -43260
$ws.Range("H139").Value = 42715
$ws.Range("J139").Value = 42715
$ws.Range("L139").Value = 42715
$ws.Range("N139").Value = -52995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 7250
$ws.Range("I19").Value = 2000
$ws.Range("J19").Value = 9000
$ws.Range("K19").Value = 2000
$ws.Range("L19").Value = 9000
$ws.Range("M19").Value = -1827
$ws.Range("N19").Value = -9346
$ws.Range("H75").Value = 13387.846
$ws.Range("I75").Value = 13836.833
$ws.Range("J75").Value = 8000
$ws.Range("K75").Value = 13836.833
$ws.Range("L75").Value = 8000
$ws.Range("M75").Value = -12900.833
$ws.Range("N75").Value = -9872
$ws.Range("H78").Value = 13387.846
$ws.Range("I78").Value = 13836.833
$ws.Range("J78").Value = 8000
$ws.Range("K78").Value = 41510.499
$ws.Range("L78").Value = 24000
$ws.Range("M78").Value = -36830.499
$ws.Range("N78").Value = -33360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1893.1538
$ws.Range("I99").Value = 1961.1
$ws.Range("J99").Value = 1666.6666
$ws.Range("K99").Value = 1961.1
$ws.Range("L99").Value = 1666.6666
$ws.Range("M99").Value = -463.0999999999999
$ws.Range("N99").Value = -4662.6666
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H111").Value = 89000
$ws.Range("J111").Value = 89000
$ws.Range("L111").Value = 89000
$ws.Range("N111").Value = -97180
$ws.Range("H126").Value = 1893.1538
$ws.Range("I126").Value = 1961.1
$ws.Range("J126").Value = 1666.6666
$ws.Range("K126").Value = 5883.299999999999
$ws.Range("L126").Value = 4999.9998
$ws.Range("M126").Value = -3413.299999999999
$ws.Range("N126").Value = -9939.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 200
$ws.Range("I51").Value = 200
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 600
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -140
$ws.Range("N51").ClearContents()
$ws.Range("H57").Value = 3862.5
$ws.Range("I57").Value = 1966.6666
$ws.Range("K57").Value = 5899.9998
$ws.Range("M57").Value = -5340.9998
$ws.Range("H58").Value = 3249.5
$ws.Range("I58").Value = 2500
$ws.Range("J58").Value = 3999
$ws.Range("K58").Value = 7500
$ws.Range("L58").Value = 11997
$ws.Range("M58").Value = -7372
$ws.Range("N58").Value = -12253
$ws.Range("H63").Value = 1599.6
$ws.Range("I63").Value = 999.5
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 2998.5
$ws.Range("L63").Value = 12000
$ws.Range("M63").Value = -2249.5
$ws.Range("N63").Value = -13498
$ws.Range("H66").Value = 1599.6
$ws.Range("I66").Value = 999.5
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 8995.5
$ws.Range("L66").Value = 36000
$ws.Range("M66").Value = -5251.5
$ws.Range("N66").Value = -43488
$ws.Range("H75").Value = 1938
$ws.Range("I75").Value = 400
$ws.Range("J75").Value = 5014
$ws.Range("K75").Value = 1200
$ws.Range("L75").Value = 15042
$ws.Range("M75").Value = -202
$ws.Range("N75").Value = -17038
$ws.Range("H78").Value = 1938
$ws.Range("I78").Value = 400
$ws.Range("J78").Value = 5014
$ws.Range("K78").Value = 3600
$ws.Range("L78").Value = 45126
$ws.Range("M78").Value = 1392
$ws.Range("N78").Value = -55110
$ws.Range("H81").Value = 129504
$ws.Range("I81").Value = 606.5
$ws.Range("J81").Value = 181063
$ws.Range("K81").Value = 1819.5
$ws.Range("L81").Value = 543189
$ws.Range("M81").Value = -696.5
$ws.Range("N81").Value = -545435
$ws.Range("H84").Value = 129504
$ws.Range("I84").Value = 606.5
$ws.Range("J84").Value = 181063
$ws.Range("K84").Value = 5458.5
$ws.Range("L84").Value = 1629567
$ws.Range("M84").Value = 157.5
$ws.Range("N84").Value = -1640799
$ws.Range("H87").Value = 33744.035
$ws.Range("I87").Value = 8498.75
$ws.Range("J87").Value = 34975.51
$ws.Range("K87").Value = 25496.25
$ws.Range("L87").Value = 104926.53
$ws.Range("M87").Value = -24248.25
$ws.Range("N87").Value = -107422.53
$ws.Range("H88").Value = 4000
$ws.Range("J88").Value = 4000
$ws.Range("L88").Value = 12000
$ws.Range("N88").Value = -12856
$ws.Range("H90").Value = 33744.035
$ws.Range("I90").Value = 8498.75
$ws.Range("J90").Value = 34975.51
$ws.Range("K90").Value = 76488.75
$ws.Range("L90").Value = 314779.59
$ws.Range("M90").Value = -70248.75
$ws.Range("N90").Value = -327259.59
$ws.Range("H91").Value = 4000
$ws.Range("J91").Value = 4000
$ws.Range("L91").Value = 12000
$ws.Range("N91").Value = -14964
$ws.Range("H131").Value = 736.3137
$ws.Range("I131").Value = 453.16666
$ws.Range("J131").Value = 890.75757
$ws.Range("K131").Value = 1359.49998
$ws.Range("L131").Value = 2672.27271
$ws.Range("M131").Value = 3680.50002
$ws.Range("N131").Value = -12752.27271

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 3333.3333
$ws.Range("J5").Value = 3333.3333
$ws.Range("L5").Value = 3333.3333
$ws.Range("N5").Value = -3557.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1427.6111
$ws.Range("I16").Value = 513.13336
$ws.Range("K16").Value = 513.13336
$ws.Range("M16").Value = -343.13336
$ws.Range("H74").Value = 42000
$ws.Range("I74").Value = 30000
$ws.Range("J74").Value = 60000
$ws.Range("K74").Value = 30000
$ws.Range("L74").Value = 60000
$ws.Range("M74").Value = -29002
$ws.Range("N74").Value = -61996
$ws.Range("H77").Value = 42000
$ws.Range("I77").Value = 30000
$ws.Range("J77").Value = 60000
$ws.Range("K77").Value = 90000
$ws.Range("L77").Value = 180000
$ws.Range("M77").Value = -85008
$ws.Range("N77").Value = -189984

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("H126").Value = 2255.913
$ws.Range("I126").Value = 2128.7
$ws.Range("J126").Value = 2353.7693
$ws.Range("K126").Value = 6386.099999999999
$ws.Range("L126").Value = 7061.3079
$ws.Range("M126").Value = -3916.099999999999
$ws.Range("N126").Value = -12001.3079
